$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'246.08"
$ws.Range("D3").Formula = "'26.14"
$ws.Range("D4").Formula = "'5.088"
$ws.Range("D5").Formula = "'0.05618"
$ws.Range("D6").Formula = "'6.479"
$ws.Range("D7").Formula = "'3.021"
$ws.Range("D8").Formula = "'0.8112"
$ws.Range("D9").Formula = "'0.8489"
$ws.Range("D10").Formula = "'0.1346"
$ws.Range("D11").Formula = "'0.03237"
$ws.Range("D12").Formula = "'0.02765"
$ws.Range("D13").Formula = "'0.09404"
$ws.Range("D14").Formula = "'0.001513"
$ws.Range("D15").Formula = "'0.0005998"
$ws.Range("D16").Formula = "'0.006082"
$ws.Range("D17").Formula = "'3.558"
$ws.Range("D20").Formula = "'0.06965"
$ws.Range("D21").Formula = "'0.1319"
$ws.Range("D22").Formula = "'3.742"
$ws.Range("D23").Formula = "'0.04697"
$ws.Range("D24").Formula = "'0.1374"
$ws.Range("D25").Formula = "'0.001246"
$ws.Range("D27").Formula = "'0.00009596"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Formula = "'0.1364"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Formula = "'0.002659"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Formula = "'0.003411"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Formula = "'0.008633"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Formula = "'0.00005291"
$ws.Range("D47").Formula = "'0.1329"
